$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 10: article (create) ---
$ws.Range("A10").Value = "article"
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial($PasteFormats)

$ws.Range("B10").Value = "https://ecos.joheee.com/googolplex/article"
[void]$ws.Hyperlinks.Add($ws.Range("B10"), "https://ecos.joheee.com/googolplex/article")
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial($PasteFormats)

$ws.Range("C10").Value = "post"
$ws.Range("C7").Copy()
$ws.Range("C10").PasteSpecial($PasteFormats)

$ws.Range("D10").Value = "{`n  `"title`": `"this is article title`",`n  `"content`": `"this is article content`"`n}"
$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial($PasteFormats)

$ws.Range("E10").Value = "{`n  `"status`": 200,`n  `"message`": `"article with title this is article title is successfully created!`",`n  `"data`": {`n    `"id`": `"68cdbf11-a175-4487-a78a-8402261ec869`",`n    `"title`": `"this is article title`",`n    `"content`": `"this is article content`"`n  }`n}"
$ws.Range("E7").Copy()
$ws.Range("E10").PasteSpecial($PasteFormats)

$ws.Range("F10").Value = "done"
$ws.Range("F7").Copy()
$ws.Range("F10").PasteSpecial($PasteFormats)

$ws.Rows.Item(10).RowHeight = 135

# --- Row 11: article (update) ---
$ws.Range("A11").Value = "article"
$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial($PasteFormats)

$ws.Range("B11").Value = "https://ecos.joheee.com/googolplex/article/68cdbf11-a175-4487-a78a-8402261ec869"
[void]$ws.Hyperlinks.Add($ws.Range("B11"), "https://ecos.joheee.com/googolplex/article/68cdbf11-a175-4487-a78a-8402261ec869")
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial($PasteFormats)

$ws.Range("C11").Value = "patch"
$ws.Range("C7").Copy()
$ws.Range("C11").PasteSpecial($PasteFormats)

$ws.Range("D11").Value = "{`n  `"title`": `"this is update article title`",`n  `"content`": `"this is update article content`"`n}"
$ws.Range("D7").Copy()
$ws.Range("D11").PasteSpecial($PasteFormats)

$ws.Range("E11").Value = "{`n  `"status`": 200,`n  `"message`": `"article is successfully updated!`",`n  `"data`": {`n    `"id`": `"68cdbf11-a175-4487-a78a-8402261ec869`",`n    `"title`": `"this is update article title`",`n    `"content`": `"this is update article content`"`n  }`n}"
$ws.Range("E7").Copy()
$ws.Range("E11").PasteSpecial($PasteFormats)

$ws.Range("F11").Value = "done"
$ws.Range("F7").Copy()
$ws.Range("F11").PasteSpecial($PasteFormats)

$ws.Rows.Item(11).RowHeight = 135

# --- Row 12: assignment (create) ---
$ws.Range("A12").Value = "assignment"
$ws.Range("A7").Copy()
$ws.Range("A12").PasteSpecial($PasteFormats)

$ws.Range("B12").Value = "https://ecos.joheee.com/googolplex/assignment"
[void]$ws.Hyperlinks.Add($ws.Range("B12"), "https://ecos.joheee.com/googolplex/assignment")
$ws.Range("B7").Copy()
$ws.Range("B12").PasteSpecial($PasteFormats)

$ws.Range("C12").Value = "post"
$ws.Range("C7").Copy()
$ws.Range("C12").PasteSpecial($PasteFormats)

$ws.Range("D12").Value = "{`n  `"title`": `"this is assignment title`",`n  `"content`": `"this is assignment content`",`n  `"due_date`": `"2024-12-31T23:59:59.999Z`"`n}"
$ws.Range("D7").Copy()
$ws.Range("D12").PasteSpecial($PasteFormats)

$ws.Rows.Item(12).RowHeight = 75

# --- selection / view ---
[void]$ws.Range("E2").Select()

